$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability values (more games simulated -> refined matrix)
$values = @{
    "B2" = 0.1428571428571428
    "C2" = 0.4285714285714285
    "P2" = 0.1428571428571428
    "S2" = 0.2857142857142857
    "P3" = 0.6666666666666666
    "S3" = 0.3333333333333333
    "P4" = 0.8
    "S4" = 0.2
    "P5" = 1.0
    "F6" = 0.07407407407407407
    "J6" = 0.2592592592592592
    "Q6" = 0.1481481481481481
    "R6" = 0.07407407407407407
    "S6" = 0.4444444444444444
    "B7" = 0.03703703703703703
    "D7" = 0.03703703703703703
    "F7" = 0.1851851851851852
    "Q7" = 0.1481481481481481
    "S7" = 0.3703703703703703
    "B8" = 0.02272727272727273
    "D8" = 0.02272727272727273
    "E8" = 0.01136363636363636
    "F8" = 0.03409090909090909
    "J8" = 0.04545454545454546
    "O8" = 0.01136363636363636
    "Q8" = 0.1704545454545454
    "R8" = 0.1363636363636364
    "S8" = 0.5454545454545454
    "B9" = 0.05555555555555555
    "F9" = 0.05555555555555555
    "J9" = 0.05555555555555555
    "Q9" = 0.1111111111111111
    "R9" = 0.05555555555555555
    "S9" = 0.6666666666666666
    "B10" = 0.01574803149606299
    "D10" = 0.01574803149606299
    "F10" = 0.07086614173228346
    "J10" = 0.07086614173228346
    "O10" = 0.02362204724409449
    "Q10" = 0.1811023622047244
    "R10" = 0.07874015748031496
    "S10" = 0.5433070866141733
    "G11" = 0.25
    "J11" = 0.04166666666666666
    "K11" = 0.2916666666666667
    "L11" = 0.3958333333333333
    "S11" = 0.02083333333333333
    "G12" = 0.7368421052631579
    "J12" = 0.1578947368421053
    "S12" = 0.1052631578947368
    "G13" = 0.4
    "J13" = 0.4
    "S13" = 0.2
    "H15" = 0.09523809523809523
    "I15" = 0.04761904761904762
    "J15" = 0.4285714285714285
    "K15" = 0.04761904761904762
    "S15" = 0.3809523809523809
    "H16" = 0.375
    "J16" = 0.125
    "K16" = 0.25
    "S16" = 0.25
    "F17" = 0.0625
    "H17" = 0.3333333333333333
    "I17" = 0.0625
    "J17" = 0.25
    "K17" = 0.0625
    "M17" = 0.04166666666666666
    "O17" = 0.04166666666666666
    "S17" = 0.1458333333333333
    "H18" = 0.2142857142857143
    "I18" = 0.07142857142857142
    "J18" = 0.3214285714285715
    "K18" = 0.1071428571428571
    "M18" = 0.03571428571428571
    "O18" = 0.07142857142857142
    "S18" = 0.1785714285714286
    "F19" = 0.01449275362318841
    "H19" = 0.3043478260869565
    "I19" = 0.06280193236714976
    "J19" = 0.3140096618357488
    "K19" = 0.1159420289855072
    "M19" = 0.00966183574879227
    "O19" = 0.04347826086956522
    "S19" = 0.1352657004830918
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
